# Logged Week 15 and simulated Week 16
# Update the cumulative Target Depth Data totals on both the OFF and DEF
# sheets (row 2 = "H" / Home splits) to reflect the newly logged/simulated
# week's passing attempts/completions.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 371
$wsOff.Range("C2").Value = 260
$wsOff.Range("D2").Value = 93
$wsOff.Range("E2").Value = 42

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 434
$wsDef.Range("C2").Value = 301
$wsDef.Range("D2").Value = 111
$wsDef.Range("E2").Value = 50
